$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 625.75
$ws.Range("I15").Value = 625.75
$ws.Range("K15").Value = 1877.25
$ws.Range("M15").Value = -1708.25

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1171.4286
$ws.Range("I40").Value = 1033.3334
$ws.Range("J40").Value = 1275
$ws.Range("K40").Value = 1033.3334
$ws.Range("L40").Value = 1275
$ws.Range("M40").Value = -858.3334
$ws.Range("N40").Value = -1625

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4090
$ws.Range("I64").Value = 3957.1428
$ws.Range("J64").Value = 4161.5386
$ws.Range("K64").Value = 3957.1428
$ws.Range("L64").Value = 4161.5386
$ws.Range("M64").Value = -3709.1428
$ws.Range("N64").Value = -4657.5386

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4090
$ws.Range("I67").Value = 3957.1428
$ws.Range("J67").Value = 4161.5386
$ws.Range("K67").Value = 3957.1428
$ws.Range("L67").Value = 4161.5386
$ws.Range("M67").Value = -3099.1428
$ws.Range("N67").Value = -5877.5386

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1122.034
$ws.Range("I137").Value = 946.4828
$ws.Range("J137").Value = 1461.4333
$ws.Range("K137").Value = 2839.4484
$ws.Range("L137").Value = 4384.2999
$ws.Range("M137").Value = -289.4484000000002
$ws.Range("N137").Value = -9484.2999

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7395.263
$ws.Range("I141").Value = 7743.5713
$ws.Range("J141").Value = 6420
$ws.Range("K141").Value = 23230.7139
$ws.Range("L141").Value = 19260
$ws.Range("M141").Value = -18050.7139
$ws.Range("N141").Value = -29620

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16707.797
$ws.Range("I32").Value = 16040.321
$ws.Range("J32").Value = 17840.484
$ws.Range("K32").Value = 16040.321
$ws.Range("L32").Value = 17840.484
$ws.Range("M32").Value = -15753.321
$ws.Range("N32").Value = -18414.484

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2908489
$ws.Range("I132").Value = 3379525.2
$ws.Range("J132").Value = 3766.3333
$ws.Range("K132").Value = 10138575.6
$ws.Range("L132").Value = 11298.9999
$ws.Range("M132").Value = -10136045.6
$ws.Range("N132").Value = -16358.9999

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4896.237
$ws.Range("I31").Value = 1044.5103
$ws.Range("J31").Value = 11886.407
$ws.Range("K31").Value = 1044.5103
$ws.Range("L31").Value = 11886.407
$ws.Range("M31").Value = -749.5102999999999
$ws.Range("N31").Value = -12476.407

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4896.237
$ws.Range("I34").Value = 1044.5103
$ws.Range("J34").Value = 11886.407
$ws.Range("K34").Value = 1044.5103
$ws.Range("L34").Value = 11886.407
$ws.Range("M34").Value = -842.5102999999999
$ws.Range("N34").Value = -12290.407

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5229.2856
$ws.Range("I62").Value = 5434.1665
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 5434.1665
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -4810.1665
$ws.Range("N62").Value = -5248

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5229.2856
$ws.Range("I65").Value = 5434.1665
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 27170.8325
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -24050.8325
$ws.Range("N65").Value = -26240

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 76924260
$ws.Range("I99").Value = 90910120
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 90910120
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -90908622
$ws.Range("N99").Value = -4996

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 76924260
$ws.Range("I126").Value = 90910120
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 272730360
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -272727890
$ws.Range("N126").Value = -10940

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2244.5
$ws.Range("I70").Value = 1116
$ws.Range("J70").Value = 3122.2222
$ws.Range("K70").Value = 3348
$ws.Range("L70").Value = 9366.6666
$ws.Range("M70").Value = -3033
$ws.Range("N70").Value = -9996.6666

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2244.5
$ws.Range("I73").Value = 1116
$ws.Range("J73").Value = 3122.2222
$ws.Range("K73").Value = 3348
$ws.Range("L73").Value = 9366.6666
$ws.Range("M73").Value = -2256
$ws.Range("N73").Value = -11550.6666

# CUL row 104
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 126374.875
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 126374.875
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 379124.625
$ws.Range("N104").Value = -384366.625
$ws.Range("M104").ClearContents()

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 16412863
$ws.Range("I132").Value = 20430196
$ws.Range("J132").Value = 8752.833000000001
$ws.Range("K132").Value = 61290588
$ws.Range("L132").Value = 26258.499
$ws.Range("M132").Value = -61288058
$ws.Range("N132").Value = -31318.499

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5439.778
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5439.778
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3867659.8
$ws.Range("I132").Value = 5434529.5
$ws.Range("J132").Value = 2714.2666
$ws.Range("K132").Value = 16303588.5
$ws.Range("L132").Value = 8142.7998
$ws.Range("M132").Value = -16301058.5
$ws.Range("N132").Value = -13202.7998

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 39150
$ws.Range("J133").Value = 39150
$ws.Range("L133").Value = 39150
$ws.Range("N133").Value = -44210

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64077.25
$ws.Range("J46").Value = 64077.25
$ws.Range("L46").Value = 64077.25
$ws.Range("N46").Value = -64539.25

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 313683.84
$ws.Range("I132").Value = 29398.889
$ws.Range("J132").Value = 1166538.6
$ws.Range("K132").Value = 88196.667
$ws.Range("L132").Value = 3499615.8
$ws.Range("M132").Value = -85666.667
$ws.Range("N132").Value = -3504675.8

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 64077.25
$ws.Range("J134").Value = 64077.25
$ws.Range("L134").Value = 192231.75
$ws.Range("N134").Value = -197301.75

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 849393.3
$ws.Range("I136").Value = 1033666.06
$ws.Range("J136").Value = 1738.8
$ws.Range("K136").Value = 3100998.18
$ws.Range("L136").Value = 5216.4
$ws.Range("M136").Value = -3098448.18
$ws.Range("N136").Value = -10316.4
